$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1976.1111
$ws.Range("J17").Value = 2110.625
$ws.Range("L17").Value = 6331.875
$ws.Range("N17").Value = -6667.875
$ws.Range("H38").Value = 177.6
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H41").Value = 7471.75
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 7471.75
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 7471.75
$ws.Range("M41").ClearContents()
$ws.Range("N41").Value = -8351.75
$ws.Range("H58").Value = 2851.8667
$ws.Range("I58").Value = 370.44446
$ws.Range("J58").Value = 3472.2222
$ws.Range("K58").Value = 1111.33338
$ws.Range("L58").Value = 10416.6666
$ws.Range("M58").Value = -961.33338
$ws.Range("N58").Value = -10716.6666
$ws.Range("H62").Value = 3069.5715
$ws.Range("I62").Value = 3501.25
$ws.Range("K62").Value = 3501.25
$ws.Range("M62").Value = -2877.25
$ws.Range("H65").Value = 3069.5715
$ws.Range("I65").Value = 3501.25
$ws.Range("K65").Value = 17506.25
$ws.Range("M65").Value = -14386.25
$ws.Range("H87").Value = 31666.666
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37496
$ws.Range("H90").Value = 31666.666
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -117480
$ws.Range("H100").Value = 3200
$ws.Range("I100").Value = 4000
$ws.Range("J100").Value = 2933.3333
$ws.Range("K100").Value = 4000
$ws.Range("L100").Value = 2933.3333
$ws.Range("M100").Value = -3459
$ws.Range("N100").Value = -4015.3333
$ws.Range("H103").Value = 349
$ws.Range("I103").Value = 349
$ws.Range("K103").Value = 1047
$ws.Range("M103").Value = -461
$ws.Range("H125").Value = 200001980
$ws.Range("I125").Value = 1000000000
$ws.Range("J125").Value = 2475
$ws.Range("K125").Value = 9000000000
$ws.Range("L125").Value = 22275
$ws.Range("M125").Value = -8999997540
$ws.Range("N125").Value = -27195
$ws.Range("H135").Value = 26320014
$ws.Range("J135").Value = 12791.272
$ws.Range("L135").Value = 115121.448
$ws.Range("N135").Value = -120191.448
$ws.Range("H138").Value = 3431.3296
$ws.Range("I138").Value = 3025.2173
$ws.Range("J138").Value = 3568.6912
$ws.Range("K138").Value = 9075.651899999999
$ws.Range("L138").Value = 10706.0736
$ws.Range("M138").Value = -3935.651899999999
$ws.Range("N138").Value = -20986.0736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 31576.803
$ws.Range("I32").Value = 30794.174
$ws.Range("K32").Value = 30794.174
$ws.Range("M32").Value = -30507.174
$ws.Range("H61").Value = 1866.8649
$ws.Range("I61").Value = 1467.2142
$ws.Range("J61").Value = 3110.2222
$ws.Range("K61").Value = 1467.2142
$ws.Range("L61").Value = 3110.2222
$ws.Range("M61").Value = -1255.2142
$ws.Range("N61").Value = -3534.2222
$ws.Range("H136").Value = 1866.8649
$ws.Range("I136").Value = 1467.2142
$ws.Range("J136").Value = 3110.2222
$ws.Range("K136").Value = 4401.642599999999
$ws.Range("L136").Value = 9330.6666
$ws.Range("M136").Value = -1851.642599999999
$ws.Range("N136").Value = -14430.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2337.4375
$ws.Range("I105").Value = 2132.3333
$ws.Range("K105").Value = 2132.3333
$ws.Range("M105").Value = -385.3332999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19612818
$ws.Range("I31").Value = 2588.0417
$ws.Range("J31").Value = 37044136
$ws.Range("K31").Value = 2588.0417
$ws.Range("L31").Value = 37044136
$ws.Range("M31").Value = -2293.0417
$ws.Range("N31").Value = -37044726
$ws.Range("H34").Value = 19612818
$ws.Range("I34").Value = 2588.0417
$ws.Range("J34").Value = 37044136
$ws.Range("K34").Value = 2588.0417
$ws.Range("L34").Value = 37044136
$ws.Range("M34").Value = -2386.0417
$ws.Range("N34").Value = -37044540
$ws.Range("H134").Value = 4234.2666
$ws.Range("I134").Value = 4311.3105
$ws.Range("K134").Value = 12933.9315
$ws.Range("M134").Value = -10398.9315
$ws.Range("H140").Value = 52926.668
$ws.Range("I140").Value = 0
$ws.Range("J140").Value = 52926.668
$ws.Range("K140").Value = 0
$ws.Range("L140").Value = 52926.668
$ws.Range("M140").ClearContents()
$ws.Range("N140").Value = -63286.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 766.75
$ws.Range("I131").Value = 370
$ws.Range("J131").Value = 801.25
$ws.Range("K131").Value = 1110
$ws.Range("L131").Value = 2403.75
$ws.Range("M131").Value = 3930
$ws.Range("N131").Value = -12483.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 99
$ws.Range("J2").Value = 102
$ws.Range("L2").Value = 102
$ws.Range("N2").Value = -328
$ws.Range("H43").Value = 1921.9
$ws.Range("H57").Value = 9700
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()
$ws.Range("H80").Value = 12000
$ws.Range("I80").Value = 9000
$ws.Range("K80").Value = 9000
$ws.Range("M80").Value = -8002
$ws.Range("H83").Value = 12000
$ws.Range("I83").Value = 9000
$ws.Range("K83").Value = 45000
$ws.Range("M83").Value = -40008

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 10000
$ws.Range("I64").Value = 10000
$ws.Range("K64").Value = 10000
$ws.Range("M64").Value = -9775
$ws.Range("H67").Value = 10000
$ws.Range("I67").Value = 10000
$ws.Range("K67").Value = 10000
$ws.Range("M67").Value = -9220
$ws.Range("H122").Value = 9000.666999999999
$ws.Range("I122").Value = 11626
$ws.Range("J122").Value = 3750
$ws.Range("K122").Value = 34878
$ws.Range("L122").Value = 11250
$ws.Range("M122").Value = -32428
$ws.Range("N122").Value = -16150

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2915.2222
$ws.Range("I122").Value = 3179.5789
$ws.Range("J122").Value = 2287.375
$ws.Range("K122").Value = 9538.736699999999
$ws.Range("L122").Value = 6862.125
$ws.Range("M122").Value = -7088.736699999999
$ws.Range("N122").Value = -11762.125
$ws.Range("H136").Value = 1804.5714
$ws.Range("I136").Value = 1721.8055
$ws.Range("J136").Value = 2301.1667
$ws.Range("K136").Value = 5165.416499999999
$ws.Range("L136").Value = 6903.500100000001
$ws.Range("M136").Value = -2615.416499999999
$ws.Range("N136").Value = -12003.5001
